$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the answer key values in column A (rows 2-21).
$ws.Range("A2").Value = "A"
$ws.Range("A3").Value = "B"
$ws.Range("A4").Value = "A"
$ws.Range("A5").Value = "B"
$ws.Range("A6").Value = "B"
$ws.Range("A7").Value = "D"
$ws.Range("A8").Value = "C"
$ws.Range("A9").Value = "B"
$ws.Range("A10").Value = "C"
$ws.Range("A11").Value = "C"
$ws.Range("A12").Value = "B"
$ws.Range("A13").Value = "D"
$ws.Range("A14").Value = "A"
$ws.Range("A15").Value = "A"
$ws.Range("A16").Value = "A"
$ws.Range("A17").Value = "C"
$ws.Range("A18").Value = "C"
$ws.Range("A19").Value = "A"
$ws.Range("A20").Value = "A"
$ws.Range("A21").Value = "B"

# Move the active selection to match the saved view state.
$ws.Range("G7").Select()
